$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("axes")

# Insert three new columns (D:F) with "*_arrow" headers / "*(%)" values,
# pushing the existing "Title" / "QAP Diagram (Intrusive)" column from D to G.
$ws.Columns("D:F").Insert()

$ws.Range("D1").Value = "A_arrow"
$ws.Range("E1").Value = "B_arrow"
$ws.Range("F1").Value = "C_arrow"

$ws.Range("D2").Value = "Q (%)"
$ws.Range("E2").Value = "A (%)"
$ws.Range("F2").Value = "P (%)"

$ws.Columns("D:F").ColumnWidth = 14

$ws.Activate()
$ws.Range("E7").Select() | Out-Null
